$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 14654
$ws1.Range("F10").Value = 15160
$ws1.Range("F11").Value = 30
$ws1.Range("F13").Value = 315
$ws1.Range("F18").Value = 182
$ws1.Range("F25").Value = 5
$ws1.Range("F37").Value = 5309

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 14654
$ws4.Range("F10").Value = 15160
$ws4.Range("F11").Value = 30
$ws4.Range("F13").Value = 315
$ws4.Range("F19").Value = 182
$ws4.Range("F26").Value = 5
$ws4.Range("F40").Value = 5309
